# Apply the change described by the diff:
#  - Clear the contents of cells in column C (rows 2-59) that currently hold 0,
#    turning them into truly empty cells (no stored value), leaving the other
#    (non-zero) values untouched.
#  - Update the active sheet's selection / scroll position: remove the
#    scrolled "topLeftCell" and move the selection from H24 to G49.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 59; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 0) {
        $cell.ClearContents()
        $cell.Style = $cell.Style
    }
}

# Update view / selection to match the saved state in the workbook.
$ws.Activate()
$ws.Range("G49").Select()
